$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)
$newRow = $tbl.Rows.Add($tbl.Rows.Item(1))

$cell1 = $tbl.Cell(1, 1)
$cell1.Range.Text = "Tipo de evento"
$cell1.Range.Font.Bold = $true
$cell1.Range.Font.BoldBi = $true

$cell2 = $tbl.Cell(1, 2)
$cell2.Range.Text = "Localización"
$cell2.Range.Font.Bold = $true
$cell2.Range.Font.BoldBi = $true

Write-Output $tbl.Rows.Count
